$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Alpha smoothing")
$ws.Activate()

# Row 20: section header
$ws.Range("A20").Value = "After Binned method"

# Row 21
$ws.Range("A21").Value = "Gr1 Muus 1000#2 pcrit 7-27-21 B.txt"
$ws.Range("B21").Value = 0.72
$ws.Range("C21").Value = 0.91
$ws.Range("D21").Value = 1.38
$ws.Range("E21").Value = 1.83
$ws.Range("F21").Value = 2.16
$ws.Range("H21").Value = 5
$ws.Range("I21").Value = 5
$ws.Range("J21").Value = 5.01
$ws.Range("K21").Value = 5.09
$ws.Range("L21").Value = 5.0999999999999996

# Row 22
$ws.Range("A22").Value = "tbocto 1800 pcrit day 7 tank 10 blank ch 3 4 8-20-21-ch1.txt"
$ws.Range("B22").Value = 0.8
$ws.Range("C22").Value = 1.03
$ws.Range("D22").Value = 1.44
$ws.Range("E22").Value = 1.79
$ws.Range("F22").Value = 1.87
$ws.Range("H22").Value = 2.6
$ws.Range("I22").Value = 2.5499999999999998
$ws.Range("J22").Value = 2.0299999999999998
$ws.Range("K22").Value = 2.35
$ws.Range("L22").Value = 2.34

# Row 23: section header
$ws.Range("A23").Value = "tbocto 1800 pcrit tank 10 3-4 blank 8-13-21.txt"

# Row 24
$ws.Range("A24").Value = "tbocto 1000 pcrit tank 3 and 4 8-11-21-ch1.txt"
$ws.Range("B24").Value = -0.81
$ws.Range("C24").Value = 1.26
$ws.Range("D24").Value = 2.0499999999999998
$ws.Range("E24").Value = 2.4900000000000002
$ws.Range("F24").Value = 2.89
$ws.Range("H24").Value = 18.5
$ws.Range("I24").Value = 18.5
$ws.Range("J24").Value = 3.83
$ws.Range("K24").Value = 3.86
$ws.Range("L24").Value = 3.83

# Row 25
$ws.Range("A25").Value = "tbocto 1000 pcrti tank 3 and 4 day 7 8-19-21-ch1.txt"
$ws.Range("B25").Value = 2.38
$ws.Range("C25").Value = 2.94
$ws.Range("D25").Value = 3.96
$ws.Range("E25").Value = 4.3899999999999997
$ws.Range("F25").Value = 4.43
$ws.Range("H25").Value = 9.5
$ws.Range("I25").Value = 8.9
$ws.Range("J25").Value = 8.8000000000000007
$ws.Range("K25").Value = 8.8000000000000007
$ws.Range("L25").Value = 8.8000000000000007

# Update selection / active cell to F22 as in diff
$ws.Range("F22").Select()

# Update window view position
$excel.ActiveWindow.Left = -34880
$excel.ActiveWindow.Top = -320
